# Auto-generated edit script: updates crypto price/volume columns (D, E)
# to match the refreshed data snapshot from the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.672.11"
$ws.Range("E2").Value = "  +4.26%  "
$ws.Range("D3").Value = "2.257.92"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E4").Value = "  +0.03%  "
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.69"
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = "  +2.71%  "
$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.10"
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = "  +4.56%  "
$ws.Range("E7").Value = "  +3.31%  "
$ws.Range("E9").Value = "  +2.20%  "
$__style = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.08"
$ws.Range("D10").Style = $__style
$ws.Range("E10").Value = "  +4.86%  "
$__style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.11"
$ws.Range("D11").Style = $__style
$ws.Range("E11").Value = "  +3.30%  "
$__style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0795"
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = "  +2.08%  "
$ws.Range("E13").Value = "  +1.47%  "
$__style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.56"
$ws.Range("D14").Style = $__style
$ws.Range("E14").Value = "  +2.81%  "
$ws.Range("D15").Value = "2.605.30"
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("E16").Value = "  +2.63%  "
$ws.Range("D17").Value = "2.255.25"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("E18").Value = "  +3.16%  "
$ws.Range("D19").Value = "41.576.03"
$__style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.38"
$ws.Range("D20").Style = $__style
$ws.Range("E20").Value = "  +10.57%  "
$ws.Range("D21").Value = "0.0₃0902"
$ws.Range("E21").Value = "  +1.92%  "
$__style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.89"
$ws.Range("D22").Style = $__style
$ws.Range("E22").Value = "  +2.52%  "
$__style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.53"
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = "  +1.65%  "
$__style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.24"
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = "  +2.21%  "
$ws.Range("E25").Value = "  +4.27%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  +5.37%  "
$__style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.01"
$ws.Range("D28").Style = $__style
$ws.Range("E28").Value = "  +4.13%  "
$__style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.47"
$ws.Range("D29").Style = $__style
$__style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.07"
$ws.Range("D30").Style = $__style
$ws.Range("E30").Value = "  -1.13%  "
$__style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.72"
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = "  +0.80%  "
$__style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "34.35"
$ws.Range("D32").Style = $__style
$ws.Range("E32").Value = "  +8.06%  "
$ws.Range("E33").Value = "  +0.02%  "
$__style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.13"
$ws.Range("D34").Style = $__style
$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0741"
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = "  +4.25%  "
$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.01"
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = "  -0.67%  "
$__style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.38"
$ws.Range("D37").Style = $__style
$ws.Range("E37").Value = "  +2.24%  "
$ws.Range("E38").Value = "  +2.75%  "
$__style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.56"
$ws.Range("D39").Style = $__style
$__style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.104"
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = "  +4.03%  "
$ws.Range("E41").Value = "  +2.88%  "
$__style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.90"
$ws.Range("D42").Style = $__style
$ws.Range("E42").Value = "  +4.26%  "
$ws.Range("D43").Value = "2.058.26"
$ws.Range("E43").Value = "  -0.73%  "
$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.56"
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = "  +2.16%  "
$__style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0277"
$ws.Range("D45").Style = $__style
$ws.Range("E45").Value = "  +2.74%  "
$__style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.13"
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = "  +2.71%  "
$ws.Range("E47").Value = "  +5.72%  "
$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.86"
$ws.Range("D48").Style = $__style
$ws.Range("E48").Value = "  +3.80%  "
$__style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "72.36"
$ws.Range("D49").Style = $__style
$ws.Range("E49").Value = "  +7.28%  "
$__style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.51"
$ws.Range("D50").Style = $__style
$ws.Range("E50").Value = "  +3.34%  "
$__style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.15"
$ws.Range("D51").Style = $__style
$ws.Range("E51").Value = "  +2.73%  "
